$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (2-11) were re-shuffled/updated (weekly refresh of the
# Fruta / hortaliza price sheet). Row 4 and row 9 are unchanged; every other row's
# date, variety, volume, price, unit, origin and per-kg fields were updated to
# reflect the new weekly values.

$ws.Range("D2").Value = 44438
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 11000
$ws.Range("M2").Value = 11500
$ws.Range("O2").Value = "Provincia del Elquí"
$ws.Range("P2").Value = 383

$ws.Range("D3").Value = 44498
$ws.Range("H3").Value = "Española"
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 8500
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8750
$ws.Range("N3").Value = "$/caja 30 unidades"
$ws.Range("P3").Value = 292
$ws.Range("Q3").Value = 30

$ws.Range("D5").Value = 44426
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 11500
$ws.Range("M5").Value = 11750
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 392

$ws.Range("D6").Value = 44426
$ws.Range("H6").Value = "Madrigal"
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 12500
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12750
$ws.Range("N6").Value = "$/caja 40 unidades"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 319
$ws.Range("Q6").Value = 40

$ws.Range("D7").Value = 44484
$ws.Range("H7").Value = "Española"
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 9500
$ws.Range("N7").Value = "$/caja 30 unidades"
$ws.Range("O7").Value = "Provincia del Elquí"
$ws.Range("P7").Value = 317
$ws.Range("Q7").Value = 30

$ws.Range("D8").Value = 44427
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 12500
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 312

$ws.Range("D10").Value = 44420
$ws.Range("H10").Value = "Madrigal"
$ws.Range("J10").Value = 800
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("N10").Value = "$/caja 40 unidades"
$ws.Range("P10").Value = 362
$ws.Range("Q10").Value = 40

$ws.Range("D11").Value = 44420
$ws.Range("J11").Value = 700
$ws.Range("K11").Value = 13000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 13500
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 338
